$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a second little table (species abundance values, for the new
# "changes in individual species abundance through time" plots) by copying
# the Response names (A2:A15) and their Mean values (C2:C15) as plain
# values - not formulas - down to B17:C30.
$ws.Range("A2:A15").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4163) | Out-Null

$ws.Range("C2:C15").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false

# The new Mean column uses a plain 2-decimal number format instead of the
# old best-fit General format.
$ws.Range("C17:C30").NumberFormat = "0.00"

# Row 22 (PC_axis1) keeps the same highlight fill used on row 7.
$ws.Range("B22").Interior.Color = 65535

# Widen column C now that it no longer uses AutoFit/BestFit.
$ws.Columns.Item(3).ColumnWidth = 14.7

# Move the active selection to reflect where the user finished editing.
$ws.Range("F25").Select() | Out-Null
